$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row for India (alphabetically between Guatemala and Italy)
$ws.Rows.Item(17).Insert()
$ws.Range("A17").Value = "India"
$ws.Range("B17").Value = "IND"
$ws.Range("C17").Value = "regional"

# Re-apply the existing sort over the expanded data range so the sortState
# metadata reflects the new extent (A2:C32)
$dataRange = $ws.Range("A2:C32")
$sortKey = $ws.Range("A2:A32")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($sortKey)
$ws.Sort.SetRange($dataRange)
$ws.Sort.Header = [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlNo
$ws.Sort.Apply()

# Update the active selection to match the saved state
[void]$ws.Range("C18").Select()
